# 50% betere pie chart
#
# Restructure the "Blad1" sheet: insert a new "totaal" row at row 3
# (summing each status column), push the existing status data down by
# one row (old rows 3-17 -> new rows 4-18), rebuild the running A-column
# counter (1..15) below the new totals row, and repoint the pie chart's
# value series at the new totals row (B3:D3) instead of the old one
# (B18:D18).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- 1. Capture the current status values (columns B:D) for rows 3-17 ---
#        (the data that needs to shift down by one row). Value2 is used
#        for reads, since plain Value reflection misbehaves on this host.
$data = @()
for ($r = 3; $r -le 17; $r++) {
    $row = @{
        B = $ws.Cells.Item($r, 2).Value2
        C = $ws.Cells.Item($r, 3).Value2
        D = $ws.Cells.Item($r, 4).Value2
    }
    $data += ,$row
}

# --- 2. Clear the old data area (A3:D18) ---
$ws.Range("A3:D18").ClearContents()

# --- 3. Write the captured B:D values back, shifted down into rows 4-18 ---
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 4
    $row = $data[$i]
    if ($null -ne $row.B) { $ws.Cells.Item($r, 2).Value = $row.B }
    if ($null -ne $row.C) { $ws.Cells.Item($r, 3).Value = $row.C }
    if ($null -ne $row.D) { $ws.Cells.Item($r, 4).Value = $row.D }
}

# --- 4. Rebuild the column-A running counter (1..15) in rows 4-18 ---
$ws.Cells.Item(4, 1).Value = 1
$ws.Range("A5").Formula = "=A4+1"
$ws.Range("A6:A18").Formula = "=1+A5"

# --- 5. New "totaal" row at row 3, summing the status columns below it ---
$ws.Cells.Item(3, 1).Value = "totaal"
$ws.Range("B3").Formula = "=SUM(B4:B18)"
$ws.Range("C3").Formula = "=SUM(C4:C18)"
$ws.Range("D3").Formula = "=SUM(D4:D18)"

# --- 6. Repoint the pie chart's value series at the new totals row ---
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Values = "=Blad1!`$B`$3:`$D`$3"

# --- 7. Update the active selection ---
$ws.Range("G1").Select()
